$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0.7437728710866437
$ws.Range("D2").Value = 0

$ws.Range("B3").Value = 0
$ws.Range("C3").Value = -0.8455971708761988
$ws.Range("D3").Value = 0

$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0.7387392772873552
$ws.Range("D4").Value = -0.7509114302606947

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0.8068560497530896
$ws.Range("D5").Value = 0

$ws.Range("B6").Value = 0
$ws.Range("C6").Value = -0.6660482183303048
$ws.Range("D6").Value = 0

$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0.8264608106165505
$ws.Range("D7").Value = 0

$ws.Range("B8").Value = 0
$ws.Range("C8").Value = -0.8730779216355401
$ws.Range("D8").Value = 0.6744148420103375

$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0.7206093807323319
$ws.Range("D9").Value = 0
